# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has completed and is in sync with en-US, refreshes the handback
# timestamps, clears the stale "handback file is not latest" error
# details, and widens the Status / Error Detail columns to fit the new
# text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Status column (zh-cn + de-de) ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-10-27 09:23:17"
$wsZhCn.Range("K3").Value = "2016-10-27 09:23:17"
# The handback-version errors are resolved; clear Error Detail but keep
# the cells as (empty) text rather than deleting them outright.
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = "Normal"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-10-27 09:23:33"
$wsDeDe.Range("K3").Value = "2016-10-27 09:23:33"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = "Normal"

# --- Column width adjustments to fit the new text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470531463623
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470531463623
